$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("P2").Value = 1.24
$ws.Range("F3").Value = 1.36
$ws.Range("G3").Value = 1.37
$ws.Range("H3").Value = 10.5
$ws.Range("I3").Value = 11.5
$ws.Range("J3").Value = 5.6
$ws.Range("P3").Value = 2.64
$ws.Range("Q3").Value = 1.57
$ws.Range("U3").Value = 2.02
$ws.Range("X3").Value = 26
$ws.Range("AA3").Value = 390
$ws.Range("AB3").Value = 11
$ws.Range("AC3").Value = 13
$ws.Range("AD3").Value = 85
$ws.Range("AE3").Value = 180
$ws.Range("AH3").Value = 25
$ws.Range("AJ3").Value = 11.5
$ws.Range("AL3").Value = 55
$ws.Range("AM3").Value = 150
$ws.Range("AN3").Value = 4.8
$ws.Range("F4").Value = 5.6
$ws.Range("G4").Value = 5.8
$ws.Range("H4").Value = 1.68
$ws.Range("I4").Value = 1.69
$ws.Range("J4").Value = 4.2
$ws.Range("K4").Value = 4.5
$ws.Range("N4").Value = 5.3
$ws.Range("O4").Value = 1.21
$ws.Range("AA4").Value = 17.5
$ws.Range("AB4").Value = 26
$ws.Range("AC4").Value = 10.5
$ws.Range("AG4").Value = 23
$ws.Range("AH4").Value = 18.5
$ws.Range("AJ4").Value = 160
$ws.Range("AK4").Value = 85
$ws.Range("M5").Value = 1.05
$ws.Range("Q5").Value = 1.67
$ws.Range("R5").Value = 1.56
$ws.Range("U5").Value = 2.6
$ws.Range("X5").Value = 22
$ws.Range("Z5").Value = 21
$ws.Range("AA5").Value = 40
$ws.Range("AB5").Value = 15.5
$ws.Range("AC5").Value = 8.6
$ws.Range("AF5").Value = 21
$ws.Range("AI5").Value = 34
$ws.Range("AJ5").Value = 40
$ws.Range("AK5").Value = 27
$ws.Range("AL5").Value = 34
$ws.Range("AN5").Value = 17.5
$ws.Range("AO5").Value = 17
$ws.Range("F6").Value = 8
$ws.Range("G6").Value = 8.199999999999999
$ws.Range("H6").Value = 1.44
$ws.Range("I6").Value = 1.45
$ws.Range("J6").Value = 5.5
$ws.Range("K6").Value = 5.6
$ws.Range("N6").Value = 6
$ws.Range("P6").Value = 2.72
$ws.Range("Q6").Value = 1.54
$ws.Range("R6").Value = 1.7
$ws.Range("X6").Value = 34
$ws.Range("Y6").Value = 12
$ws.Range("Z6").Value = 10.5
$ws.Range("AA6").Value = 13
$ws.Range("AB6").Value = 36
$ws.Range("AC6").Value = 13
$ws.Range("AH6").Value = 23
$ws.Range("AI6").Value = 28
$ws.Range("AL6").Value = 1000
$ws.Range("AN6").Value = 120
$ws.Range("AO6").Value = 5
$ws.Range("F7").Value = 2.3
$ws.Range("G7").Value = 2.34
$ws.Range("H7").Value = 3.3
$ws.Range("I7").Value = 3.4
$ws.Range("M7").Value = 1.06
$ws.Range("X7").Value = 18
$ws.Range("AA7").Value = 70
$ws.Range("AB7").Value = 12.5
$ws.Range("AC7").Value = 8.4
$ws.Range("AM7").Value = 80
$ws.Range("AO7").Value = 30
$ws.Range("H8").Value = 5.1
$ws.Range("I8").Value = 5.4
$ws.Range("N8").Value = 5.6
$ws.Range("P8").Value = 2.6
$ws.Range("U8").Value = 2.48
$ws.Range("AA8").Value = 130
$ws.Range("AC8").Value = 10.5
$ws.Range("AD8").Value = 21
$ws.Range("AE8").Value = 60
$ws.Range("AJ8").Value = 18.5
$ws.Range("AN8").Value = 7.2
$ws.Range("F9").Value = 1.42
$ws.Range("G9").Value = 1.43
$ws.Range("H9").Value = 8.4
$ws.Range("I9").Value = 9
$ws.Range("J9").Value = 5.4
$ws.Range("P9").Value = 2.7
$ws.Range("R9").Value = 1.68
$ws.Range("S9").Value = 2.38
$ws.Range("U9").Value = 2.14
$ws.Range("Z9").Value = 80
$ws.Range("AI9").Value = 90
$ws.Range("AM9").Value = 95
$ws.Range("AN9").Value = 5
$ws.Range("I10").Value = 9.800000000000001
$ws.Range("O10").Value = 1.11
$ws.Range("T10").Value = 1.66
$ws.Range("AJ10").Value = 13
$ws.Range("AL10").Value = 26
$ws.Range("I11").Value = 3.35
$ws.Range("J11").Value = 3.6
$ws.Range("Z11").Value = 24
$ws.Range("AB11").Value = 11
$ws.Range("AE11").Value = 38
$ws.Range("AF11").Value = 16
$ws.Range("AH11").Value = 16.5
$ws.Range("AI11").Value = 46
$ws.Range("AK11").Value = 25
$ws.Range("AL11").Value = 44
$ws.Range("AN11").Value = 17.5
$ws.Range("AO11").Value = 32
$ws.Range("F12").Value = 1.95
$ws.Range("G12").Value = 2.1
$ws.Range("H12").Value = 3.8
$ws.Range("I12").Value = 4.2
$ws.Range("J12").Value = 3.7
$ws.Range("K12").Value = 4.1
$ws.Range("P12").Value = 2.08
$ws.Range("Q12").Value = 1.8
$ws.Range("F13").Value = 1.81
$ws.Range("G13").Value = 2
$ws.Range("H13").Value = 4.7
$ws.Range("I13").Value = 6.6
$ws.Range("K13").Value = 4.3
$ws.Range("P13").Value = 1.66
$ws.Range("Q13").Value = 2.06
